# "filter tube rows in multiform files"
#
# - The old "2006" sheet is renamed to "Feuille inutile" and its content is
#   replaced with a short "Coucou / les / zamiiiiiiiiiiiiis" placeholder.
# - A brand-new "2006" sheet is (re)created right after "Feuille inutile",
#   carrying over the original "2006" data (values + formatting) untouched.
# - The "2005" sheet gets four new header/title rows inserted above its
#   existing data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Split "2006" into "Feuille inutile" (renamed, cleared) + a fresh
#    "2006" sheet that keeps the original data/style.
# ---------------------------------------------------------------------
$old2006 = $wb.Worksheets.Item("2006")

$new2006 = $wb.Worksheets.Add($null, $old2006)
$old2006.Range("A1:B6").Copy($new2006.Range("A1"))
$new2006.Name = "__2006_tmp__"

$old2006.Cells.Clear()
$old2006.Name = "Feuille inutile"

$fiRow1 = $old2006.Range("A1:B1")
$fiRow1.Font.Name = "Arial"
$fiRow1.Font.ThemeColor = 1
$old2006.Range("A1").Value = "Coucou"
$old2006.Range("B1").Value = "les"

$fiRow2 = $old2006.Range("A2")
$fiRow2.Font.Name = "Arial"
$fiRow2.Font.ThemeColor = 1
$old2006.Range("A2").Value = "zamiiiiiiiiiiiiis"

$new2006.Name = "2006"

# ---------------------------------------------------------------------
# 2) "2005": insert 3 new rows above the existing table (title, subtitle,
#    section banner) and turn the former header row (now row 4, holding
#    the old "id de sujet" / "niveau max" labels) into the new bold
#    "Sujet (code ou autre)" / "Niveau max" header.
# ---------------------------------------------------------------------
$ws2005 = $wb.Worksheets.Item("2005")
$ws2005.Rows("1:3").Insert()

$row1 = $ws2005.Range("A1:B1")
$row1.Font.Name = "Arial"
$row1.Font.ThemeColor = 1
$ws2005.Range("A1").Value = "Parcours thématique_RGPD & Données personnelles"

$row2 = $ws2005.Range("A2:B2")
$row2.Font.Name = "Arial"
$row2.Font.ThemeColor = 1
$ws2005.Range("A2").Value = "Création : Avril 2022"

$row3 = $ws2005.Range("A3:B3")
$row3.Font.Name = "Arial"
$row3.Font.Color = 16777215
$row3.Interior.Color = 6710886
$row3.Interior.PatternColor = 6710886
$ws2005.Range("A3").Value = "1) capage"

$row4 = $ws2005.Range("A4:B4")
$row4.Font.Name = "Arial"
$row4.Font.ThemeColor = 1
$row4.Font.Bold = $true
$ws2005.Range("A4").Value = "Sujet (code ou autre)"
$ws2005.Range("B4").Value = "Niveau max"
